$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All source cells are text (inlineStr) in the workbook, including ones that
# look numeric (e.g. '11.40', '0.0000280'); force Text number format first so
# Excel's COM layer doesn't silently coerce them to Number and normalize/trim them.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.391.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.422.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.74"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.415.87"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.92"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "688.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.970.18"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.431.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.421.46"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.895"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "570.10"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.67"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.99"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.587.95"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.332"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0416"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.44"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.76"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.31%  "
